# Update spreadsheets for starting costs
# Fixes the "econ_sartingcost_duration_*" typo -> "econ_startingcost_duration_*"
# and adds the missing "econ_startingcost_duration_smearacf" row (to match the
# unit/inflection/startingcost_duration/saturation pattern used by the other
# economic programs).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# --- Row 104: fix typo on the vaccination starting-cost-duration row ---
$ws.Range("A104").Value = "econ_startingcost_duration_vaccination"

# --- Row 108: fix typo on the ipt starting-cost-duration row ---
$ws.Range("A108").Value = "econ_startingcost_duration_ipt"

# --- Row 112: fix typo on the xpert starting-cost-duration row ---
$ws.Range("A112").Value = "econ_startingcost_duration_xpert"

# --- Row 116: fix typo on the treatment_support starting-cost-duration row ---
$ws.Range("A116").Value = "econ_startingcost_duration_treatment_support"

# Insert a new row above the current "econ_saturation_smearacf" row (row 120)
# so the smearacf block gains its own starting-cost-duration row, matching
# the layout already used by vaccination/ipt/xpert/treatment_support/xpertacf.
$ws.Range("A120:E120").Insert(-4121)

# --- New row 120: starting cost duration for smearacf ---
$ws.Range("A120").Value = "econ_startingcost_duration_smearacf"
$ws.Range("B120").Value = 1

# --- Row 124 (previously 123, shifted down by the insert above): fix typo on
#     the xpertacf starting-cost-duration row ---
$ws.Range("A124").Value = "econ_startingcost_duration_xpertacf"

$wb.Save()
